$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a "Shortage of items" report. Rows 7..24 each hold one
# item, sorted alphabetically by item name (column C). A new item,
# "ANTINAL 220MG/5ML 60ML SUSP.", needs to be inserted between
# "AMBEZIM-G 30 F.C. TABS." (row 8) and "AVEROZOLID 600MG 10 F.C.TABLETS"
# (row 9), i.e. a whole new row must be inserted at row 9, pushing
# everything below (including the totals row and the footer row) down
# by one.

# Insert a new row above row 9; this shifts rows 9-26 down to 10-27 and
# adjusts the existing merged cells automatically, like Excel's own
# "Insert Sheet Rows" command run from row 9's context menu.
$ws.Rows.Item(9).Insert()

# Copy formatting (styles, merges, row height) from the row above
# (row 8, "AMBEZIM-G...") onto the freshly inserted row 9, then fix the
# row height to match this report's "…75pt" data rows.
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(9).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(9).RowHeight = 24.75

# Fill in the new row's data. Column layout (matching the other item
# rows): A = sequence number, C = item name, H = current balance,
# L = reorder limit, N = price, P = sale price, Q = transaction count.
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 3).Value = "ANTINAL 220MG/5ML 60ML SUSP."
$ws.Cells.Item(9, 8).Value = "7:0"

# Columns L and P carry genuine numeric display formats (custom / 0.00)
# on this template, yet every data row stores its value as literal text
# (matches the source report generator). Force text storage here too,
# then restore the original number format so the visual style is
# unchanged.
$cellL = $ws.Cells.Item(9, 12)
$origFmtL = $cellL.NumberFormat
$cellL.NumberFormat = "@"
$cellL.Value = "1"
$cellL.NumberFormat = $origFmtL

$ws.Cells.Item(9, 14).Value = "24.00"

$cellP = $ws.Cells.Item(9, 16)
$origFmtP = $cellP.NumberFormat
$cellP.NumberFormat = "@"
$cellP.Value = "24.0000"
$cellP.NumberFormat = $origFmtP

$ws.Cells.Item(9, 17).Value = "1:0"

# Renumber the sequence column (A) for every row that shifted down one
# position (rows 10-25 now hold what used to be rows 9-24).
for ($r = 10; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# Update the running total (a plain number, not a formula) to include
# the newly added item's price (24.00): 1526.03 + 24.00 = 1550.03.
$ws.Cells.Item(26, 14).Value = 1550.03
